$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 8 new rows after row 2 to make room for the additional match records
$ws.Rows("3:10").Insert()

# Force columns G:K (runs, balls, numberOf4, numberOf6, sr) to be stored as text
# so numeric-looking values keep their exact textual representation
# (e.g. "100.00", "87.50", "270.00") instead of being coerced to numbers.
$ws.Range("G2:K10").NumberFormat = "@"

$data = @(
    @(" Nov 2 2020", " Abu Dhabi", "Capitals won by 6 wickets (with 6 balls remaining)", "Royal Challengers Bangalore", "Delhi Capitals", "Shivam Dube ", "17", "11", "2", "1", "154.54"),
    @(" Nov 6 2020", " Abu Dhabi", "Sunrisers won by 6 wickets (with 2 balls remaining)", "Royal Challengers Bangalore", "Sunrisers Hyderabad", "Shivam Dube ", "8", "13", "0", "0", "61.53"),
    @(" Oct 15 2020", " Sharjah", "Kings XI won by 8 wickets", "Royal Challengers Bangalore", "Kings XI Punjab", "Shivam Dube ", "23", "19", "0", "2", "121.05"),
    @(" Sep 24 2020", " Dubai (DSC)", "Kings XI won by 97 runs", "Royal Challengers Bangalore", "Kings XI Punjab", "Shivam Dube ", "12", "12", "0", "1", "100.00"),
    @(" Sep 21 2020", " Dubai (DSC)", "RCB won by 10 runs", "Royal Challengers Bangalore", "Sunrisers Hyderabad", "Shivam Dube ", "7", "8", "0", "0", "87.50"),
    @(" Oct 28 2020", " Abu Dhabi", "Mumbai won by 5 wickets (with 5 balls remaining)", "Royal Challengers Bangalore", "Mumbai Indians", "Shivam Dube ", "2", "6", "0", "0", "33.33"),
    @(" Sep 28 2020", " Dubai (DSC)", "Match tied (RCB won the one-over eliminator)", "Royal Challengers Bangalore", "Mumbai Indians", "Shivam Dube ", "27", "10", "1", "3", "270.00"),
    @(" Oct 10 2020", " Dubai (DSC)", "RCB won by 37 runs", "Royal Challengers Bangalore", "Chennai Super Kings", "Shivam Dube ", "22", "14", "2", "1", "157.14"),
    @(" Oct 5 2020", " Dubai (DSC)", "Capitals won by 59 runs", "Royal Challengers Bangalore", "Delhi Capitals", "Shivam Dube ", "11", "12", "0", "1", "91.66")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $colNum = $j + 1
        $ws.Cells.Item($rowNum, $colNum).Value = $rowVals[$j]
    }
}
